$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Dinesh"
$ws.Range("B2").Value = "Mundhe"

$ws.Range("B3").Select()
